$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F21").Value = -2.349648272185756
$ws.Range("F23").Value = -2.368109808649089
$ws.Range("F24").Value = -2.376783434023615
$ws.Range("F25").Value = -2.3843324541613
$ws.Range("F29").Value = -2.420585867905181
$ws.Range("F30").Value = -2.423667425663171
$ws.Range("F32").Value = -2.440482808455025
$ws.Range("F33").Value = -2.446419973622497
$ws.Range("F34").Value = -2.455168905171035
$ws.Range("F35").Value = -2.458128281827924
$ws.Range("F36").Value = -2.46853290937509
$ws.Range("F40").Value = -2.509411304177774
$ws.Range("F41").Value = -2.506214442923099
$ws.Range("F44").Value = -2.537315501662524
$ws.Range("F46").Value = -2.546090119230279
$ws.Range("F47").Value = -2.562883810883627
$ws.Range("F48").Value = -2.57031090676471
$ws.Range("F49").Value = -2.571573488553983
$ws.Range("F50").Value = -2.586965032841244
$ws.Range("F51").Value = -2.593768830777631
$ws.Range("F55").Value = -2.61751562196456
$ws.Range("F56").Value = -2.623240480802455
$ws.Range("F57").Value = -2.632279431653752
$ws.Range("F58").Value = -2.638813278362771
$ws.Range("F59").Value = -2.645774093033165
$ws.Range("F60").Value = -2.651746091127429
$ws.Range("F61").Value = -2.658248875226303
$ws.Range("F62").Value = -2.670596732759923
$ws.Range("F63").Value = -2.673259824171872
$ws.Range("F64").Value = -2.683071514953063
$ws.Range("F66").Value = -2.688058046580476
$ws.Range("F67").Value = -2.69616750031848
$ws.Range("F70").Value = -2.709263485683898
$ws.Range("F73").Value = -2.717168498618808
$ws.Range("F74").Value = -2.718938340858696
$ws.Range("F82").Value = -2.598172641155528
$ws.Range("F83").Value = -2.641955507780406
$ws.Range("F84").Value = -2.688554140518475
$ws.Range("F85").Value = -2.738855370998668
$ws.Range("F89").Value = -2.929395910502844
$ws.Range("F94").Value = -2.726356326001451
$ws.Range("F95").Value = -2.779164846439915
$ws.Range("F96").Value = -2.830180630500621
$ws.Range("F98").Value = -2.936914993038444
$ws.Range("F99").Value = -2.987766745675359
$ws.Range("F100").Value = -3.034284751588802
$ws.Range("F103").Value = -2.703177631847244
$ws.Range("F104").Value = -2.747709196668912
$ws.Range("F105").Value = -2.796015904887811
$ws.Range("F106").Value = -2.84483660773738
$ws.Range("F108").Value = -2.950737429942385
$ws.Range("F110").Value = -3.059301430581334
$ws.Range("F111").Value = -3.106812457811352
$ws.Range("F114").Value = -2.742679176857312
$ws.Range("F116").Value = -2.835072072776722
$ws.Range("F117").Value = -2.884703791762313
$ws.Range("F118").Value = -2.936708992417926
$ws.Range("F119").Value = -2.988455583040241
$ws.Range("F120").Value = -3.045411136167934
$ws.Range("F121").Value = -3.101456316365162
$ws.Range("F125").Value = -2.777364456245955
$ws.Range("F127").Value = -2.86797705816827
$ws.Range("F128").Value = -2.915786814460953
$ws.Range("F129").Value = -2.96711231326316
$ws.Range("F130").Value = -3.019435401203694
$ws.Range("F132").Value = -3.133124833905483
$ws.Range("F133").Value = -3.189694422762355
$ws.Range("F136").Value = -2.80008832012893
$ws.Range("F138").Value = -2.889127887752382
$ws.Range("F139").Value = -2.936938551462407
$ws.Range("F140").Value = -2.987646941108987
$ws.Range("F141").Value = -3.040498161706536
$ws.Range("F142").Value = -3.097727232878441
$ws.Range("F143").Value = -3.15482555998565
$ws.Range("F144").Value = -3.21450067587261
$ws.Range("F148").Value = -2.862908399661157
$ws.Range("F149").Value = -2.907506750742893
$ws.Range("F150").Value = -2.954789591614572
$ws.Range("F151").Value = -3.005082869807285
$ws.Range("F152").Value = -3.057484266476415
$ws.Range("F153").Value = -3.11481117355711
$ws.Range("F154").Value = -3.173553799697333
$ws.Range("F158").Value = -2.835395942130353
$ws.Range("F159").Value = -2.876282456082528
$ws.Range("F161").Value = -2.96743005541098
$ws.Range("F162").Value = -3.017679903172262
$ws.Range("F163").Value = -3.071744347277336
$ws.Range("F165").Value = -3.188009369529822
$ws.Range("F166").Value = -3.250429172508638
$ws.Range("F171").Value = -2.933392568684186
$ws.Range("F172").Value = -2.979521849811723
$ws.Range("F174").Value = -3.083039235306934
$ws.Range("F175").Value = -3.14076851060391
$ws.Range("F177").Value = -3.264051237610057
$ws.Range("F180").Value = -2.85877499928102
$ws.Range("F181").Value = -2.89921218433041
$ws.Range("F182").Value = -2.941993504784392
$ws.Range("F183").Value = -2.988778607546368
$ws.Range("F184").Value = -3.039054842342528
$ws.Range("F185").Value = -3.09361778735938
$ws.Range("F186").Value = -3.150711593252427
$ws.Range("F188").Value = -3.2752812884611
$ws.Range("F191").Value = -2.868334518258368
$ws.Range("F192").Value = -2.9081195526928
$ws.Range("F193").Value = -2.95090094949013
$ws.Range("F194").Value = -2.997213616483758
$ws.Range("F195").Value = -3.047499457761634
$ws.Range("F203").Value = -2.915432953831704
$ws.Range("F205").Value = -3.004681459623628
$ws.Range("F206").Value = -3.054986595682447
$ws.Range("F207").Value = -3.109213788794329
$ws.Range("F208").Value = -3.166848793693372
$ws.Range("F209").Value = -3.228459968776472
$ws.Range("F210").Value = -3.293695001232593
